$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Update header role labels (row 1) in sharedStrings used by the sheet
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "Procurement Sourcing Leader"
$ws.Range("D1").Value = "S3L "
$ws.Range("E1").Value = "GSC" + [char]10 + "Manager"
$ws.Range("G1").Value = "Regional Category Manager"
$ws.Range("J1").Value = "NAM" + [char]10 + "Regional Category Director"
$ws.Range("K1").Value = "VP Category Manager"

# ---------------------------------------------------------------------
# Add a new cell H4 that duplicates the value already present in G4
# ---------------------------------------------------------------------
$ws.Range("H4").Value = $ws.Range("G4").Value2

# ---------------------------------------------------------------------
# Row height adjustments
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 53.3
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(9).RowHeight = 27.2

# ---------------------------------------------------------------------
# Default column width tweak
# ---------------------------------------------------------------------
$ws.StandardWidth = 8.5390625

# ---------------------------------------------------------------------
# View state: scroll back to A1, zoom out to 63%, move the selection
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("E18").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 63
